$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" "63.304.12"
Set-TextValue "E2" "  +1.19%  "
Set-TextValue "D3" "2.457.68"
Set-TextValue "E3" "  +0.94%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "572.87"
Set-TextValue "E5" "  +1.04%  "
Set-TextValue "D6" "146.81"
Set-TextValue "E6" "  +0.85%  "
Set-TextValue "E7" "  +0.07%  "
Set-TextValue "D8" "0.538"
Set-TextValue "E8" "  +1.08%  "
Set-TextValue "D9" "2.454.51"
Set-TextValue "E9" "  +0.76%  "
Set-TextValue "D10" "0.112"
Set-TextValue "E10" "  +0.97%  "
Set-TextValue "D11" "0.157"
Set-TextValue "E11" "  +1.38%  "
Set-TextValue "D12" "5.29"
Set-TextValue "E12" "  -0.36%  "
Set-TextValue "D13" "0.357"
Set-TextValue "E13" "  +0.47%  "
Set-TextValue "D14" "27.09"
Set-TextValue "E14" "  +0.89%  "
Set-TextValue "D15" "0.0000180"
Set-TextValue "E15" "  +0.19%  "
Set-TextValue "D16" "2.900.26"
Set-TextValue "E16" "  +1.39%  "
Set-TextValue "D17" "63.210.80"
Set-TextValue "E17" "  +1.28%  "
Set-TextValue "D18" "2.450.66"
Set-TextValue "E18" "  +0.92%  "
Set-TextValue "D19" "11.30"
Set-TextValue "E19" "  +0.55%  "
Set-TextValue "D20" "7.34"
Set-TextValue "E20" "  +5.05%  "
Set-TextValue "D21" "329.27"
Set-TextValue "E21" "  +1.65%  "
Set-TextValue "D22" "4.22"
Set-TextValue "E22" "  +1.09%  "
Set-TextValue "E23" "  +13.50%  "
Set-TextValue "D24" "1.00"
Set-TextValue "E24" "  +0.19%  "
Set-TextValue "D25" "65.76"
Set-TextValue "E25" "  -2.17%  "
Set-TextValue "D26" "619.42"
Set-TextValue "E26" "  +4.85%  "
Set-TextValue "D27" "8.99"
Set-TextValue "E27" "  +5.06%  "
Set-TextValue "D28" "0.0000103"
Set-TextValue "E28" "  +2.67%  "
Set-TextValue "D29" "2.563.92"
Set-TextValue "E29" "  +0.54%  "
Set-TextValue "D30" "1.51"
Set-TextValue "E30" "  +4.68%  "
Set-TextValue "D32" "8.27"
Set-TextValue "E32" "  -2.11%  "
Set-TextValue "D33" "1.90"
Set-TextValue "E33" "  +1.31%  "
Set-TextValue "E34" "  -3.04%  "
Set-TextValue "D35" "5.20"
Set-TextValue "E35" "  +7.30%  "
Set-TextValue "E36" "  +1.89%  "
Set-TextValue "E37" "  +0.12%  "
Set-TextValue "E38" "  -0.42%  "
Set-TextValue "D39" "5.44"
Set-TextValue "E39" "  +1.76%  "
Set-TextValue "D40" "18.88"
Set-TextValue "E40" "  +0.71%  "
Set-TextValue "D41" "146.91"
Set-TextValue "E41" "  -0.59%  "
Set-TextValue "D42" "1.79"
Set-TextValue "E42" "  -2.12%  "
Set-TextValue "D43" "2.61"
Set-TextValue "E43" "  +6.68%  "
Set-TextValue "E44" "  -0.03%  "
Set-TextValue "D45" "41.86"
Set-TextValue "E45" "  +0.62%  "
Set-TextValue "D46" "148.79"
Set-TextValue "D47" "3.77"
Set-TextValue "E47" "  +2.77%  "
Set-TextValue "D48" "21.26"
Set-TextValue "D49" "0.0535"
Set-TextValue "E49" "  -0.05%  "
Set-TextValue "D50" "0.603"
Set-TextValue "E50" "  +0.32%  "
Set-TextValue "D51" "0.0233"
Set-TextValue "E51" "  +0.83%  "
